# "documented parts list"
#
# 1. Reword the column-C header ("Sequence (including restriction sites, cut
#    sites)" -> "sequence (with restriction sites, cut sites, etc.)").
# 2. Fix a typo'd BsaI/BsmBI cut-site motif baked into many of the part
#    sequences in column C: "...TCTTCTG..." -> "...TCGTCTG...".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-word the sequence column header in row 1 -----------------------
$ws.Range("C1").Value = "sequence (with restriction sites, cut sites, etc.)"

# --- 2. Correct the cut-site typo in every affected sequence in column C --
$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("C$r")
    $val = $cell.Value2
    if ($val -ne $null -and $val.ToString().Contains("TCTTCTG")) {
        $cell.Value = $val.ToString().Replace("TCTTCTG", "TCGTCTG")
    }
}
